$d = $word.ActiveDocument

$d.Content.Find.Execute("2+61=63", $true, $true, $false, $false, $false, $true, 1, $false, "43+43=86", 2) | Out-Null
$d.Content.Find.Execute("55-39=16", $true, $true, $false, $false, $false, $true, 1, $false, "22+53=75", 2) | Out-Null
$d.Content.Find.Execute("97-61=36", $true, $true, $false, $false, $false, $true, 1, $false, "39-21=18", 2) | Out-Null
$d.Content.Find.Execute("2+94=96", $true, $true, $false, $false, $false, $true, 1, $false, "75-54=21", 2) | Out-Null
$d.Content.Find.Execute("28+50=78", $true, $true, $false, $false, $false, $true, 1, $false, "31-21=10", 2) | Out-Null
$d.Content.Find.Execute("29+25=54", $true, $true, $false, $false, $false, $true, 1, $false, "13+36=49", 2) | Out-Null
$d.Content.Find.Execute("37+48=85", $true, $true, $false, $false, $false, $true, 1, $false, "42-36=6", 2) | Out-Null
$d.Content.Find.Execute("26+64=90", $true, $true, $false, $false, $false, $true, 1, $false, "33+14=47", 2) | Out-Null
$d.Content.Find.Execute("92-38=54", $true, $true, $false, $false, $false, $true, 1, $false, "67-58=9", 2) | Out-Null
$d.Content.Find.Execute("98-20=78", $true, $true, $false, $false, $false, $true, 1, $false, "43+33=76", 2) | Out-Null
$d.Content.Find.Execute("84-69=15", $true, $true, $false, $false, $false, $true, 1, $false, "10+71=81", 2) | Out-Null
$d.Content.Find.Execute("63-62=1", $true, $true, $false, $false, $false, $true, 1, $false, "77-15=62", 2) | Out-Null
$d.Content.Find.Execute("97-96=1", $true, $true, $false, $false, $false, $true, 1, $false, "25+13=38", 2) | Out-Null
$d.Content.Find.Execute("12+28=40", $true, $true, $false, $false, $false, $true, 1, $false, "89-39=50", 2) | Out-Null
$d.Content.Find.Execute("26+18=44", $true, $true, $false, $false, $false, $true, 1, $false, "52+7=59", 2) | Out-Null
$d.Content.Find.Execute("23+15=38", $true, $true, $false, $false, $false, $true, 1, $false, "6+88=94", 2) | Out-Null
$d.Content.Find.Execute("19+50=69", $true, $true, $false, $false, $false, $true, 1, $false, "61-43=18", 2) | Out-Null
$d.Content.Find.Execute("47-45=2", $true, $true, $false, $false, $false, $true, 1, $false, "10+70=80", 2) | Out-Null
$d.Content.Find.Execute("15-1=14", $true, $true, $false, $false, $false, $true, 1, $false, "24+56=80", 2) | Out-Null
$d.Content.Find.Execute("54-2=52", $true, $true, $false, $false, $false, $true, 1, $false, "25-2=23", 2) | Out-Null
$d.Content.Find.Execute("20+26=46", $true, $true, $false, $false, $false, $true, 1, $false, "52-6=46", 2) | Out-Null
$d.Content.Find.Execute("40+23=63", $true, $true, $false, $false, $false, $true, 1, $false, "21-8=13", 2) | Out-Null
$d.Content.Find.Execute("51-16=35", $true, $true, $false, $false, $false, $true, 1, $false, "21+74=95", 2) | Out-Null
$d.Content.Find.Execute("18+11=29", $true, $true, $false, $false, $false, $true, 1, $false, "43+36=79", 2) | Out-Null
$d.Content.Find.Execute("46-43=3", $true, $true, $false, $false, $false, $true, 1, $false, "16-13=3", 2) | Out-Null
$d.Content.Find.Execute("33+28=61", $true, $true, $false, $false, $false, $true, 1, $false, "35+64=99", 2) | Out-Null
$d.Content.Find.Execute("35-25=10", $true, $true, $false, $false, $false, $true, 1, $false, "76-54=22", 2) | Out-Null
$d.Content.Find.Execute("81+9=90", $true, $true, $false, $false, $false, $true, 1, $false, "25+35=60", 2) | Out-Null
$d.Content.Find.Execute("18-3=15", $true, $true, $false, $false, $false, $true, 1, $false, "25+72=97", 2) | Out-Null
$d.Content.Find.Execute("30+35=65", $true, $true, $false, $false, $false, $true, 1, $false, "5+65=70", 2) | Out-Null
$d.Content.Find.Execute("93+1=94", $true, $true, $false, $false, $false, $true, 1, $false, "45+21=66", 2) | Out-Null
$d.Content.Find.Execute("92-29=63", $true, $true, $false, $false, $false, $true, 1, $false, "0+31=31", 2) | Out-Null
$d.Content.Find.Execute("51+0=51", $true, $true, $false, $false, $false, $true, 1, $false, "80-49=31", 2) | Out-Null
$d.Content.Find.Execute("61-44=17", $true, $true, $false, $false, $false, $true, 1, $false, "65-41=24", 2) | Out-Null
$d.Content.Find.Execute("32+1=33", $true, $true, $false, $false, $false, $true, 1, $false, "7+1=8", 2) | Out-Null
$d.Content.Find.Execute("53-0=53", $true, $true, $false, $false, $false, $true, 1, $false, "87-85=2", 2) | Out-Null
$d.Content.Find.Execute("67-61=6", $true, $true, $false, $false, $false, $true, 1, $false, "72-28=44", 2) | Out-Null
$d.Content.Find.Execute("9+8=17", $true, $true, $false, $false, $false, $true, 1, $false, "38+16=54", 2) | Out-Null
$d.Content.Find.Execute("94-31=63", $true, $true, $false, $false, $false, $true, 1, $false, "39+29=68", 2) | Out-Null
$d.Content.Find.Execute("96-17=79", $true, $true, $false, $false, $false, $true, 1, $false, "81-69=12", 2) | Out-Null
$d.Content.Find.Execute("46+37=83", $true, $true, $false, $false, $false, $true, 1, $false, "66+31=97", 2) | Out-Null
$d.Content.Find.Execute("8-4=4", $true, $true, $false, $false, $false, $true, 1, $false, "80-61=19", 2) | Out-Null
$d.Content.Find.Execute("5+51=56", $true, $true, $false, $false, $false, $true, 1, $false, "21+37=58", 2) | Out-Null
$d.Content.Find.Execute("95-71=24", $true, $true, $false, $false, $false, $true, 1, $false, "57-25=32", 2) | Out-Null
$d.Content.Find.Execute("49-21=28", $true, $true, $false, $false, $false, $true, 1, $false, "23+74=97", 2) | Out-Null
$d.Content.Find.Execute("85+11=96", $true, $true, $false, $false, $false, $true, 1, $false, "88-9=79", 2) | Out-Null
$d.Content.Find.Execute("64+33=97", $true, $true, $false, $false, $false, $true, 1, $false, "48+22=70", 2) | Out-Null
$d.Content.Find.Execute("3+91=94", $true, $true, $false, $false, $false, $true, 1, $false, "20+68=88", 2) | Out-Null
$d.Content.Find.Execute("19+3=22", $true, $true, $false, $false, $false, $true, 1, $false, "74-7=67", 2) | Out-Null
$d.Content.Find.Execute("46+18=64", $true, $true, $false, $false, $false, $true, 1, $false, "64-56=8", 2) | Out-Null
$d.Content.Find.Execute("95-45=50", $true, $true, $false, $false, $false, $true, 1, $false, "40-30=10", 2) | Out-Null
$d.Content.Find.Execute("14+12=26", $true, $true, $false, $false, $false, $true, 1, $false, "62-58=4", 2) | Out-Null
$d.Content.Find.Execute("86-57=29", $true, $true, $false, $false, $false, $true, 1, $false, "69-61=8", 2) | Out-Null
$d.Content.Find.Execute("22+60=82", $true, $true, $false, $false, $false, $true, 1, $false, "6+2=8", 2) | Out-Null
$d.Content.Find.Execute("43+45=88", $true, $true, $false, $false, $false, $true, 1, $false, "90-74=16", 2) | Out-Null
$d.Content.Find.Execute("83+8=91", $true, $true, $false, $false, $false, $true, 1, $false, "89+7=96", 2) | Out-Null
$d.Content.Find.Execute("61+16=77", $true, $true, $false, $false, $false, $true, 1, $false, "2+20=22", 2) | Out-Null
$d.Content.Find.Execute("65+17=82", $true, $true, $false, $false, $false, $true, 1, $false, "77+2=79", 2) | Out-Null
$d.Content.Find.Execute("85-35=50", $true, $true, $false, $false, $false, $true, 1, $false, "16-6=10", 2) | Out-Null
$d.Content.Find.Execute("80+4=84", $true, $true, $false, $false, $false, $true, 1, $false, "94-50=44", 2) | Out-Null
$d.Content.Find.Execute("33+2=35", $true, $true, $false, $false, $false, $true, 1, $false, "30+10=40", 2) | Out-Null
$d.Content.Find.Execute("36-0=36", $true, $true, $false, $false, $false, $true, 1, $false, "81-38=43", 2) | Out-Null
$d.Content.Find.Execute("56-24=32", $true, $true, $false, $false, $false, $true, 1, $false, "12+38=50", 2) | Out-Null
$d.Content.Find.Execute("78-9=69", $true, $true, $false, $false, $false, $true, 1, $false, "2+97=99", 2) | Out-Null
$d.Content.Find.Execute("2+80=82", $true, $true, $false, $false, $false, $true, 1, $false, "25-22=3", 2) | Out-Null
$d.Content.Find.Execute("23+42=65", $true, $true, $false, $false, $false, $true, 1, $false, "48+16=64", 2) | Out-Null
$d.Content.Find.Execute("56-16=40", $true, $true, $false, $false, $false, $true, 1, $false, "3+76=79", 2) | Out-Null
$d.Content.Find.Execute("5+60=65", $true, $true, $false, $false, $false, $true, 1, $false, "9-5=4", 2) | Out-Null
$d.Content.Find.Execute("89-12=77", $true, $true, $false, $false, $false, $true, 1, $false, "22+41=63", 2) | Out-Null
$d.Content.Find.Execute("54+15=69", $true, $true, $false, $false, $false, $true, 1, $false, "39+2=41", 2) | Out-Null
$d.Content.Find.Execute("60-27=33", $true, $true, $false, $false, $false, $true, 1, $false, "76+16=92", 2) | Out-Null
$d.Content.Find.Execute("20+53=73", $true, $true, $false, $false, $false, $true, 1, $false, "41-12=29", 2) | Out-Null
$d.Content.Find.Execute("9+70=79", $true, $true, $false, $false, $false, $true, 1, $false, "92+6=98", 2) | Out-Null
$d.Content.Find.Execute("14+59=73", $true, $true, $false, $false, $false, $true, 1, $false, "83-73=10", 2) | Out-Null
$d.Content.Find.Execute("46+9=55", $true, $true, $false, $false, $false, $true, 1, $false, "96-67=29", 2) | Out-Null
$d.Content.Find.Execute("62+2=64", $true, $true, $false, $false, $false, $true, 1, $false, "62-30=32", 2) | Out-Null
$d.Content.Find.Execute("29+38=67", $true, $true, $false, $false, $false, $true, 1, $false, "74+12=86", 2) | Out-Null
$d.Content.Find.Execute("43-15=28", $true, $true, $false, $false, $false, $true, 1, $false, "96-60=36", 2) | Out-Null
$d.Content.Find.Execute("42+49=91", $true, $true, $false, $false, $false, $true, 1, $false, "31-16=15", 2) | Out-Null
$d.Content.Find.Execute("5+47=52", $true, $true, $false, $false, $false, $true, 1, $false, "19+45=64", 2) | Out-Null
$d.Content.Find.Execute("66+10=76", $true, $true, $false, $false, $false, $true, 1, $false, "77+7=84", 2) | Out-Null
$d.Content.Find.Execute("23+7=30", $true, $true, $false, $false, $false, $true, 1, $false, "2+65=67", 2) | Out-Null
$d.Content.Find.Execute("12+77=89", $true, $true, $false, $false, $false, $true, 1, $false, "53-17=36", 2) | Out-Null
$d.Content.Find.Execute("65-18=47", $true, $true, $false, $false, $false, $true, 1, $false, "14+80=94", 2) | Out-Null
$d.Content.Find.Execute("62-28=34", $true, $true, $false, $false, $false, $true, 1, $false, "76-6=70", 2) | Out-Null
$d.Content.Find.Execute("67-0=67", $true, $true, $false, $false, $false, $true, 1, $false, "40+10=50", 2) | Out-Null
$d.Content.Find.Execute("11+17=28", $true, $true, $false, $false, $false, $true, 1, $false, "0+22=22", 2) | Out-Null
$d.Content.Find.Execute("96-69=27", $true, $true, $false, $false, $false, $true, 1, $false, "41+35=76", 2) | Out-Null
$d.Content.Find.Execute("48-48=0", $true, $true, $false, $false, $false, $true, 1, $false, "12+22=34", 2) | Out-Null
$d.Content.Find.Execute("2+78=80", $true, $true, $false, $false, $false, $true, 1, $false, "84-35=49", 2) | Out-Null
$d.Content.Find.Execute("70-28=42", $true, $true, $false, $false, $false, $true, 1, $false, "84-2=82", 2) | Out-Null
$d.Content.Find.Execute("14+0=14", $true, $true, $false, $false, $false, $true, 1, $false, "63-49=14", 2) | Out-Null
$d.Content.Find.Execute("41+56=97", $true, $true, $false, $false, $false, $true, 1, $false, "98-7=91", 2) | Out-Null
$d.Content.Find.Execute("62-17=45", $true, $true, $false, $false, $false, $true, 1, $false, "81-13=68", 2) | Out-Null
$d.Content.Find.Execute("92-82=10", $true, $true, $false, $false, $false, $true, 1, $false, "20+43=63", 2) | Out-Null
$d.Content.Find.Execute("75+2=77", $true, $true, $false, $false, $false, $true, 1, $false, "5+21=26", 2) | Out-Null
$d.Content.Find.Execute("91-49=42", $true, $true, $false, $false, $false, $true, 1, $false, "17+78=95", 2) | Out-Null
$d.Content.Find.Execute("21+41=62", $true, $true, $false, $false, $false, $true, 1, $false, "70+22=92", 2) | Out-Null
$d.Content.Find.Execute("86-75=11", $true, $true, $false, $false, $false, $true, 1, $false, "39+20=59", 2) | Out-Null
$d.Content.Find.Execute("63-53=10", $true, $true, $false, $false, $false, $true, 1, $false, "19+29=48", 2) | Out-Null
